# T1_algo_benchmark.xlsx - update benchmark numbers (re-run of the scoring
# script produced slightly different metrics) and refresh a handful of
# heat-map cell colors plus the sheet view/column-width cosmetics that
# come along with re-exporting the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- updated benchmark values (PDB / ArchiveII columns, rows 4-7) -------
$ws.Range("B4").Value = 0.9003357831124881
$ws.Range("C4").Value = 0.9121918550193481
$ws.Range("D4").Value = 0.8917096772664029
$ws.Range("E4").Value = 0.548525044703939
$ws.Range("F4").Value = 0.6023897386797379
$ws.Range("G4").Value = 0.5716313096962714

$ws.Range("B5").Value = 0.8839677876605528
$ws.Range("C5").Value = 0.9105430677875705
$ws.Range("D5").Value = 0.8814311239920871
$ws.Range("E5").Value = 0.5734367102821211
$ws.Range("F5").Value = 0.6362392054665977
$ws.Range("G5").Value = 0.6003153627108263

$ws.Range("B6").Value = 0.9145469687416055
$ws.Range("C6").Value = 0.9286118077613434
$ws.Range("D6").Value = 0.9048169259454163
$ws.Range("E6").Value = 0.7276246130475307
$ws.Range("F6").Value = 0.7602157620995333
$ws.Range("G6").Value = 0.7398378918339839

$ws.Range("B7").Value = 0.8102598320746768
$ws.Range("C7").Value = 0.9650308865834685
$ws.Range("D7").Value = 0.8697750695805618
$ws.Range("E7").Value = 0.8252173220350971
$ws.Range("F7").Value = 0.8865941821421769
$ws.Range("G7").Value = 0.8515154081751982

# --- heat-map fill colors that shifted along with the new values --------
$ws.Range("F4").Interior.Color = 8302892   # 2CB17E
$ws.Range("G4").Interior.Color = 8563493   # 25AB82
$ws.Range("C5").Interior.Color = 2089418   # CAE11F
$ws.Range("F5").Interior.Color = 7911479   # 37B878
$ws.Range("B6").Interior.Color = 1958349   # CDE11D

# --- cosmetic sheet-view / column-width refresh from the re-export ------
$ws.Columns("A").ColumnWidth = 14.25
$ws.Range("B1:M1").EntireColumn.ColumnWidth = 7

$ws.Range("A1:M8").Select()
$excel.ActiveWindow.Zoom = 179
